$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New aircraft data rows (A: Name, B: EW, C: MTOW)
$data = @(
    @("C-130",   34686,  70305),
    @("C-17",    128140, 265352),
    @("C-5",     169640, 379655),
    @("Convair", 32579,  74843),
    @("Hughes",  122500, 181500)
)

$row = 11
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}

# Match the style used by rows B/C in rows 8-10 (numFmt thousands separator, right aligned)
$ws.Range("B11:C15").NumberFormat = "#,##0"
$ws.Range("B11:C15").HorizontalAlignment = -4152

# Update selection to D14 (as in the diff)
$ws.Range("D14").Select()
